$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = "meta-instanceID"
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = "uuid"
$ws.Range("D42").Value = "instanceID"

$ws.Range("G7").Select()
